# Generate Report for Handback
# Regenerates the handback-status workbook's file identifiers:
#   4745c70d-b4ea-450f-9f38-8cbb53bf8401 -> 52ca86cf-fe3f-49d9-8e54-fafdcc507556
#   bbb21bc0-aa5a-45b0-b82a-aa27382279ef -> ffffdd7731fd-164b-49d2-a5db-953eff5d4133
#   xlf content hashes (d9c94df..., 0d18fcb6...) -> e1fb3f4bd307738e79e1e5b8f297a97cb62c9848
# and refreshes the handoff/handback timestamps, while leaving the
# underlying hyperlink targets (URLs / rels) untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: A2/A3 are links to the per-language markdown source
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$ov = @(
    @{ Cell = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/55b966345826d5706157efb80ddaaebf55eb7711/e2e/4745c70d-b4ea-450f-9f38-8cbb53bf8401.md"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.md" },
    @{ Cell = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/55b966345826d5706157efb80ddaaebf55eb7711/e2e/bbb21bc0-aa5a-45b0-b82a-aa27382279ef.md"; Display = "ffffdd7731fd-164b-49d2-a5db-953eff5d4133.md" }
)

$wsOverview.Hyperlinks.Delete()
foreach ($item in $ov) {
    $wsOverview.Hyperlinks.Add($wsOverview.Range($item.Cell), $item.Address, $null, $null, $item.Display)
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zh = @(
    @{ Cell = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/55b966345826d5706157efb80ddaaebf55eb7711/e2e/4745c70d-b4ea-450f-9f38-8cbb53bf8401.md"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.md" },
    @{ Cell = "D2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d41107b1d9cf90ea52f08ae6432144e6c8e32b05/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4745c70d-b4ea-450f-9f38-8cbb53bf8401.d9c94df57918df55139f8e10704db13e11879e7b.zh-cn.xlf"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.e1fb3f4bd307738e79e1e5b8f297a97cb62c9848.zh-cn.xlf" },
    @{ Cell = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4808a864d158967ec812ad70a4784734dc88f2f8/e2e/4745c70d-b4ea-450f-9f38-8cbb53bf8401.md"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.md" },
    @{ Cell = "G2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/88015da7684797d289867d5c61279c55a44af975/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4745c70d-b4ea-450f-9f38-8cbb53bf8401.d9c94df57918df55139f8e10704db13e11879e7b.zh-cn.xlf"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.e1fb3f4bd307738e79e1e5b8f297a97cb62c9848.zh-cn.xlf" },
    @{ Cell = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/55b966345826d5706157efb80ddaaebf55eb7711/e2e/bbb21bc0-aa5a-45b0-b82a-aa27382279ef.md"; Display = "ffffdd7731fd-164b-49d2-a5db-953eff5d4133.md" },
    @{ Cell = "D3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d41107b1d9cf90ea52f08ae6432144e6c8e32b05/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bbb21bc0-aa5a-45b0-b82a-aa27382279ef.0d18fcb63c5252685ddbdabe43580147162c1e42.zh-cn.xlf"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.e1fb3f4bd307738e79e1e5b8f297a97cb62c9848.zh-cn.xlf" },
    @{ Cell = "F3"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4808a864d158967ec812ad70a4784734dc88f2f8/e2e/bbb21bc0-aa5a-45b0-b82a-aa27382279ef.md"; Display = "ffffdd7731fd-164b-49d2-a5db-953eff5d4133.md" },
    @{ Cell = "G3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/88015da7684797d289867d5c61279c55a44af975/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bbb21bc0-aa5a-45b0-b82a-aa27382279ef.0d18fcb63c5252685ddbdabe43580147162c1e42.zh-cn.xlf"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.e1fb3f4bd307738e79e1e5b8f297a97cb62c9848.zh-cn.xlf" }
)

$wsZh.Hyperlinks.Delete()
foreach ($item in $zh) {
    $wsZh.Hyperlinks.Add($wsZh.Range($item.Cell), $item.Address, $null, $null, $item.Display)
}

# Handoff / handback timestamps (non-hyperlinked cells)
$wsZh.Range("E2").Value = "2016-03-22 13:13:05"
$wsZh.Range("H2").Value = "2016-03-22 13:13:24"
$wsZh.Range("E3").Value = "2016-03-22 13:13:05"
$wsZh.Range("H3").Value = "2016-03-22 13:13:24"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$de = @(
    @{ Cell = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/55b966345826d5706157efb80ddaaebf55eb7711/e2e/4745c70d-b4ea-450f-9f38-8cbb53bf8401.md"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.md" },
    @{ Cell = "D2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff46116709763493ab360fa103858ca20e2892c8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4745c70d-b4ea-450f-9f38-8cbb53bf8401.d9c94df57918df55139f8e10704db13e11879e7b.de-de.xlf"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.e1fb3f4bd307738e79e1e5b8f297a97cb62c9848.de-de.xlf" },
    @{ Cell = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6c7e4ce1eca612bce158a8cfc46089040559930a/e2e/4745c70d-b4ea-450f-9f38-8cbb53bf8401.md"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.md" },
    @{ Cell = "G2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2d0ba627ad2b4c59227e1bdf87635531c19e0812/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4745c70d-b4ea-450f-9f38-8cbb53bf8401.d9c94df57918df55139f8e10704db13e11879e7b.de-de.xlf"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.e1fb3f4bd307738e79e1e5b8f297a97cb62c9848.de-de.xlf" },
    @{ Cell = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/55b966345826d5706157efb80ddaaebf55eb7711/e2e/bbb21bc0-aa5a-45b0-b82a-aa27382279ef.md"; Display = "ffffdd7731fd-164b-49d2-a5db-953eff5d4133.md" },
    @{ Cell = "D3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff46116709763493ab360fa103858ca20e2892c8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bbb21bc0-aa5a-45b0-b82a-aa27382279ef.0d18fcb63c5252685ddbdabe43580147162c1e42.de-de.xlf"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.e1fb3f4bd307738e79e1e5b8f297a97cb62c9848.de-de.xlf" },
    @{ Cell = "F3"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/6c7e4ce1eca612bce158a8cfc46089040559930a/e2e/bbb21bc0-aa5a-45b0-b82a-aa27382279ef.md"; Display = "ffffdd7731fd-164b-49d2-a5db-953eff5d4133.md" },
    @{ Cell = "G3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2d0ba627ad2b4c59227e1bdf87635531c19e0812/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bbb21bc0-aa5a-45b0-b82a-aa27382279ef.0d18fcb63c5252685ddbdabe43580147162c1e42.de-de.xlf"; Display = "52ca86cf-fe3f-49d9-8e54-fafdcc507556.e1fb3f4bd307738e79e1e5b8f297a97cb62c9848.de-de.xlf" }
)

$wsDe.Hyperlinks.Delete()
foreach ($item in $de) {
    $wsDe.Hyperlinks.Add($wsDe.Range($item.Cell), $item.Address, $null, $null, $item.Display)
}

# Handoff / handback timestamps (non-hyperlinked cells)
$wsDe.Range("E2").Value = "2016-03-22 13:13:09"
$wsDe.Range("H2").Value = "2016-03-22 13:13:31"
$wsDe.Range("E3").Value = "2016-03-22 13:13:09"
$wsDe.Range("H3").Value = "2016-03-22 13:13:31"
